# repull data, push all data, mean calculation
# Update column F (dSF) with repulled values for the affected rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    3  = 4
    4  = 10
    7  = 2
    8  = 2
    10 = 0
    11 = -4
    15 = -1
    16 = 4
    19 = -3
    22 = -1
    26 = -2
    33 = -4
    34 = -2
    35 = 0
    40 = -4
    43 = 4
    46 = 0
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
